$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.087.34"
$ws.Range("E2").Value = "  -2.22%  "
$ws.Range("D3").Value = "2.172.23"
$ws.Range("E3").Value = "  -2.14%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'237.03"
$ws.Range("E5").Value = "  -2.38%  "
$ws.Range("D6").Value = "'0.614"
$ws.Range("E6").Value = "  -0.85%  "
$ws.Range("D7").Value = "'70.09"
$ws.Range("E7").Value = "  -4.94%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").Value = "'0.578"
$ws.Range("E9").Value = "  -6.22%  "
$ws.Range("D10").Value = "'40.12"
$ws.Range("E10").Value = "  -8.20%  "
$ws.Range("D11").Value = "'0.0929"
$ws.Range("E11").Value = "  -3.20%  "
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").Value = "'0.101"
$ws.Range("E12").Value = "  -2.69%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'6.76"
$ws.Range("E13").Value = "  -5.13%  "
$ws.Range("D14").Value = "2.497.24"
$ws.Range("E14").Value = "  -2.15%  "
$ws.Range("D15").Value = "'13.91"
$ws.Range("E15").Value = "  -2.33%  "
$ws.Range("D16").Value = "'0.809"
$ws.Range("E16").Value = "  -4.25%  "
$ws.Range("D17").Value = "2.161.83"
$ws.Range("E17").Value = "  -3.38%  "
$ws.Range("D18").Value = "40.925.33"
$ws.Range("E18").Value = "  -2.33%  "
$ws.Range("D20").Value = "'70.40"
$ws.Range("E20").Value = "  -2.75%  "
$ws.Range("D21").Value = "'5.93"
$ws.Range("E21").Value = "  -4.20%  "
$ws.Range("D22").Value = "'10.05"
$ws.Range("E22").Value = "  -3.82%  "
$ws.Range("D23").Value = "'225.62"
$ws.Range("E24").Value = "  -7.77%  "
$ws.Range("E25").Value = "  -0.02%  "
$ws.Range("D26").Value = "'10.88"
$ws.Range("E26").Value = "  -5.68%  "
$ws.Range("D27").Value = "'3.52"
$ws.Range("E27").Value = "  -2.05%  "
$ws.Range("D28").Value = "'2.20"
$ws.Range("E28").Value = "  -3.40%  "
$ws.Range("D29").Value = "'2.18"
$ws.Range("E29").Value = "  -1.27%  "
$ws.Range("E30").Value = "  +0.17%  "
$ws.Range("D31").Value = "'19.92"
$ws.Range("E31").Value = "  -3.24%  "
$ws.Range("D32").Value = "'30.86"
$ws.Range("E32").Value = "  +5.37%  "
$ws.Range("D33").Value = "'0.0769"
$ws.Range("E33").Value = "  -3.73%  "
$ws.Range("E34").Value = "  -8.98%  "
$ws.Range("E35").Value = "  -3.18%  "
$ws.Range("E36").Value = "  -8.40%  "
$ws.Range("E37").Value = "  -3.30%  "
$ws.Range("D38").Value = "'0.0285"
$ws.Range("E38").Value = "  -5.16%  "
$ws.Range("E39").Value = "  -4.67%  "
$ws.Range("E40").Value = "  -3.12%  "
$ws.Range("D41").Value = "'5.42"
$ws.Range("E41").Value = "  -4.40%  "
$ws.Range("D42").Value = "'60.24"
$ws.Range("E42").Value = "  -7.90%  "
$ws.Range("D43").Value = "'0.190"
$ws.Range("E43").Value = "  -4.93%  "
$ws.Range("E44").Value = "  -4.71%  "
$ws.Range("D45").Value = "'0.0969"
$ws.Range("E45").Value = "  -4.17%  "
$ws.Range("D46").Value = "'98.28"
$ws.Range("E46").Value = "  -6.03%  "
$ws.Range("E47").Value = "  -2.68%  "
$ws.Range("E48").Value = "  -2.83%  "
$ws.Range("E49").Value = "  -7.38%  "
$ws.Range("E50").Value = "  -3.01%  "
$ws.Range("D51").Value = "2.373.98"
$ws.Range("E51").Value = "  -2.17%  "
